# Add Location data to "All_SANs" and "SANs in Stock" (Darwin_Items / Darwin_Timestamps)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# All_SANs: add "Location" header column (D) and three new inventory rows
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("All_SANs")

$ws.Range("D1").Value = "Location"
$ws.Range("D1").Style = "Normal"

$ws.Range("A121").Value = "SAN126998"
$ws.Range("B121").Value = "Laptop 840 G10"
$ws.Range("C121").Value = "2024-11-17 14:54:43"
$ws.Range("A121:C121").Style = "Normal"

$ws.Range("A122").Value = "SAN126999"
$ws.Range("B122").Value = "Laptop 840 G10"
$ws.Range("C122").Value = "2024-11-17 14:54:47"
$ws.Range("A122:C122").Style = "Normal"

$ws.Range("A123").Value = "SAN126985"
$ws.Range("B123").Value = "Desktop Mini G9"
$ws.Range("C123").Value = "2024-11-17 14:56:48"
$ws.Range("D123").Value = "Darwin"
$ws.Range("A123:D123").Style = "Normal"

# ---------------------------------------------------------------------------
# Darwin_Items: update LastCount/NewCount for the items affected by the
# new Darwin stock movements
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Darwin_Items")

$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 2

# ---------------------------------------------------------------------------
# Darwin_Timestamps: clear the trailing empty "SAN #" placeholder cells on
# the existing rows, then log the new stock movements
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Darwin_Timestamps")

$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()

$ws.Range("A10").Value = "2024-11-17 14:54:43"
$ws.Range("B10").Value = "Laptop 840 G10"
$ws.Range("C10").Value = "add"
$ws.Range("D10").Value = "SAN126998"

$ws.Range("A11").Value = "2024-11-17 14:54:47"
$ws.Range("B11").Value = "Laptop 840 G10"
$ws.Range("C11").Value = "add"
$ws.Range("D11").Value = "SAN126999"

$ws.Range("A12").Value = "2024-11-17 14:54:47"
$ws.Range("B12").Value = "Laptop 840 G10"
$ws.Range("C12").Value = "add 2"

$ws.Range("A13").Value = "2024-11-17 14:56:43"
$ws.Range("B13").Value = "Desktop Mini G9"
$ws.Range("C13").Value = "add"
$ws.Range("D13").Value = "SAN126987"

$ws.Range("A14").Value = "2024-11-17 14:56:48"
$ws.Range("B14").Value = "Desktop Mini G9"
$ws.Range("C14").Value = "add"
$ws.Range("D14").Value = "SAN126985"

$ws.Range("A15").Value = "2024-11-17 14:56:49"
$ws.Range("B15").Value = "Desktop Mini G9"
$ws.Range("C15").Value = "add 2"

$ws.Range("A16").Value = "2024-11-17 14:58:12"
$ws.Range("B16").Value = "Desktop Mini G9"
$ws.Range("C16").Value = "subtract"
$ws.Range("D16").Value = "SAN126987"

$ws.Range("A17").Value = "2024-11-17 14:58:12"
$ws.Range("B17").Value = "Desktop Mini G9"
$ws.Range("C17").Value = "subtract 1"

# ---------------------------------------------------------------------------
# BR_Timestamps / L17_Timestamps: clear the other trailing empty "SAN #"
# placeholder cells left over from earlier rows
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BR_Timestamps")
$ws.Range("D20").ClearContents()
$ws.Range("D21").ClearContents()

$ws = $wb.Worksheets.Item("L17_Timestamps")
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
